$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D. This shifts the existing D..H data
# (Compensation %, fiberAeff_1, fiberAeff_2, fiberAlphadB_1, fiberAlphadB_2)
# one column to the right (now E..I), making room for a new "Span (km)"
# column of data at D.
$ws.Columns("D:D").Insert()

# Copy the header formatting (bold font + border) from an existing header
# cell onto the newly inserted header cell, then set its text.
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 4).PasteSpecial(-4122)
$ws.Cells.Item(1, 4).Value = "Span (km)"

# Fill in the Span (km) value (30) for every data row (rows 2-14).
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 4).Value = 30
}
